$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A12:J12").Copy()
$ws.Range("A13:J13").PasteSpecial()

$ws.Range("A13").Value = "噩梦花瓣"
$ws.Range("B13").Value = "petals_evil"
$ws.Range("C13").Value = "blue"
$ws.Range("D13").Value = "credit_coins"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 10
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = "normal"
$ws.Range("I13").Value = "official"
$ws.Range("J13").Value = "official"

$ws.Range("I19").Select()
